$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Bdnf"
$ws.Cells.Item(2,3).Value = "Ntrk2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.03254066666666667
$ws.Cells.Item(2,8).Value = 0.097622
$ws.Cells.Item(2,9).Value = 0.02633076257175775
$ws.Cells.Item(2,10).Value = 0.03898291886711458
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 1.098765
$ws.Cells.Item(2,14).Value = 3.296295
$ws.Cells.Item(2,15).Value = 0.04219121913039277
$ws.Cells.Item(2,16).Value = 0.04777771124863077
$ws.Cells.Item(2,17).Value = 0.03575454561
$ws.Cells.Item(2,18).Value = 0.32179091049
$ws.Cells.Item(2,19).Value = 0.001110926973535375
$ws.Cells.Item(2,20).Value = 0.001862514641261801

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Bdnf"
$ws.Cells.Item(3,3).Value = "Ntrk2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.03254066666666667
$ws.Cells.Item(3,8).Value = 0.097622
$ws.Cells.Item(3,9).Value = 0.02633076257175775
$ws.Cells.Item(3,10).Value = 0.03898291886711458
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 15.79937666666667
$ws.Cells.Item(3,14).Value = 47.39813
$ws.Cells.Item(3,15).Value = 0.6066765532820465
$ws.Cells.Item(3,16).Value = 0.6870059169052114
$ws.Cells.Item(3,17).Value = 0.5141222496511112
$ws.Cells.Item(3,18).Value = 4.62710024686
$ws.Cells.Item(3,19).Value = 0.0159742562823219
$ws.Cells.Item(3,20).Value = 0.02678149591994352

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Bdnf"
$ws.Cells.Item(4,3).Value = "Ntrk2"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.03254066666666667
$ws.Cells.Item(4,8).Value = 0.097622
$ws.Cells.Item(4,9).Value = 0.02633076257175775
$ws.Cells.Item(4,10).Value = 0.03898291886711458
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.009167
$ws.Cells.Item(4,14).Value = 0.027501
$ws.Cells.Item(4,15).Value = 0.000352001479632415
$ws.Cells.Item(4,16).Value = 0.0003986096017039115
$ws.Cells.Item(4,17).Value = 0.0002983002913333334
$ws.Cells.Item(4,18).Value = 0.002684702622
$ws.Cells.Item(4,19).Value = 0.00000926846738510854
$ws.Cells.Item(4,20).Value = 0.00001553896576287644

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Bdnf"
$ws.Cells.Item(5,3).Value = "Ntrk2"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.03254066666666667
$ws.Cells.Item(5,8).Value = 0.097622
$ws.Cells.Item(5,9).Value = 0.02633076257175775
$ws.Cells.Item(5,10).Value = 0.03898291886711458
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 9.1351955
$ws.Cells.Item(5,14).Value = 18.270391
$ws.Cells.Item(5,15).Value = 0.3507802261079284
$ws.Cells.Item(5,16).Value = 0.264817762244454
$ws.Cells.Item(5,17).Value = 0.2972653517003334
$ws.Cells.Item(5,18).Value = 1.783592110202
$ws.Cells.Item(5,19).Value = 0.00923631084851536
$ws.Cells.Item(5,20).Value = 0.01032336934014639

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Bdnf"
$ws.Cells.Item(6,3).Value = "Ntrk2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.2033015
$ws.Cells.Item(6,8).Value = 2.406603
$ws.Cells.Item(6,9).Value = 0.9736692374282422
$ws.Cells.Item(6,10).Value = 0.9610170811328854
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.098765
$ws.Cells.Item(6,14).Value = 3.296295
$ws.Cells.Item(6,15).Value = 0.04219121913039277
$ws.Cells.Item(6,16).Value = 0.04777771124863077
$ws.Cells.Item(6,17).Value = 1.3221455726475
$ws.Cells.Item(6,18).Value = 7.932873435885
$ws.Cells.Item(6,19).Value = 0.04108029215685739
$ws.Cells.Item(6,20).Value = 0.04591519660736897

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Bdnf"
$ws.Cells.Item(7,3).Value = "Ntrk2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.2033015
$ws.Cells.Item(7,8).Value = 2.406603
$ws.Cells.Item(7,9).Value = 0.9736692374282422
$ws.Cells.Item(7,10).Value = 0.9610170811328854
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 15.79937666666667
$ws.Cells.Item(7,14).Value = 47.39813
$ws.Cells.Item(7,15).Value = 0.6066765532820465
$ws.Cells.Item(7,16).Value = 0.6870059169052114
$ws.Cells.Item(7,17).Value = 19.011413642065
$ws.Cells.Item(7,18).Value = 114.06848185239
$ws.Cells.Item(7,19).Value = 0.5907022969997245
$ws.Cells.Item(7,20).Value = 0.6602244209852679

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Bdnf"
$ws.Cells.Item(8,3).Value = "Ntrk2"
$ws.Cells.Item(8,4).Value = "M1"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.2033015
$ws.Cells.Item(8,8).Value = 2.406603
$ws.Cells.Item(8,9).Value = 0.9736692374282422
$ws.Cells.Item(8,10).Value = 0.9610170811328854
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.009167
$ws.Cells.Item(8,14).Value = 0.027501
$ws.Cells.Item(8,15).Value = 0.000352001479632415
$ws.Cells.Item(8,16).Value = 0.0003986096017039115
$ws.Cells.Item(8,17).Value = 0.0110306648505
$ws.Cells.Item(8,18).Value = 0.06618398910300001
$ws.Cells.Item(8,19).Value = 0.0003427330122473064
$ws.Cells.Item(8,20).Value = 0.000383070635941035

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Bdnf"
$ws.Cells.Item(9,3).Value = "Ntrk2"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.2033015
$ws.Cells.Item(9,8).Value = 2.406603
$ws.Cells.Item(9,9).Value = 0.9736692374282422
$ws.Cells.Item(9,10).Value = 0.9610170811328854
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 9.1351955
$ws.Cells.Item(9,14).Value = 18.270391
$ws.Cells.Item(9,15).Value = 0.3507802261079284
$ws.Cells.Item(9,16).Value = 0.264817762244454
$ws.Cells.Item(9,17).Value = 10.99239444794325
$ws.Cells.Item(9,18).Value = 43.969577791773
$ws.Cells.Item(9,19).Value = 0.341543915259413
$ws.Cells.Item(9,20).Value = 0.2544943929043076
